# Revert "version con revision de estilo":
#  - Active tab moves from "CUADERNO DEL PROFESOR" (4th sheet) back to "GUION" (1st sheet)
#  - Selection on "CUADERNO DEL PROFESOR" resets from C14 (B14:C14) to A2 (A2:C24)
#  - The "APROVECHADO" (si/no) column on "CUADERNO DEL PROFESOR" reverts to its
#    pre-style-review values for rows 2-21 (header in row 1)

$wb = $excel.ActiveWorkbook

$ws_guion   = $wb.Worksheets.Item("GUION")
$ws_profesor = $wb.Worksheets.Item("CUADERNO DEL PROFESOR")

# --- Revert the APROVECHADO (yes/no) values on "CUADERNO DEL PROFESOR" ---
$ws_profesor.Activate()

$ws_profesor.Range("C2").Value  = "no"
$ws_profesor.Range("C3").Value  = "no"
$ws_profesor.Range("C4").Value  = "no"
$ws_profesor.Range("C6").Value  = "no"
$ws_profesor.Range("C7").Value  = "no"
$ws_profesor.Range("C8").Value  = "no"
$ws_profesor.Range("C9").Value  = "no"
$ws_profesor.Range("C10").Value = "no"
$ws_profesor.Range("C11").Value = "no"
$ws_profesor.Range("C12").Value = "no"
$ws_profesor.Range("C13").Value = "no"
$ws_profesor.Range("C14").Value = "sí"
$ws_profesor.Range("C16").Value = "no"
$ws_profesor.Range("C17").Value = "no"
$ws_profesor.Range("C19").Value = "no"
$ws_profesor.Range("C20").Value = "no"
$ws_profesor.Range("C21").Value = "no"

# --- Reset the remembered selection on "CUADERNO DEL PROFESOR" ---
$ws_profesor.Range("A2:C24").Select() | Out-Null

# --- Move the active tab back to "GUION" ---
$ws_guion.Activate()
$ws_guion.Range("B2").Select() | Out-Null
